$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Row 1 header updates
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 updates - C2 and E2 are cleared entirely
$ws.Range("B2").Value = 15.607610807449943
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 21.379660433761398
$ws.Range("E2").ClearContents()

# Row 3 updates
$ws.Range("B3").Value = 12.126629948255992
$ws.Range("C3").Value = -12.284955612774002
$ws.Range("D3").Value = 19.679804590674184
$ws.Range("E3").Value = -10.224154572232669

# Update selection to match new used-data range
$ws.Range("B1:E3").Select()
